# This edit removes specific companies / duplicate-scope rows from the
# three data sheets of the workbook, mirroring a data-cleanup commit
# ("Calculate and plot scores").
#
# Sheet "target_data":
#   - Allegro's separate Scope-3 row is dropped (its Scope-1+2 row stays).
#   - Orange's separate S1+S3 row is dropped (its Scope-1+2 row stays).
#   - JSW, Grupa Kety, Eurocash and CitiBank Handlowy are removed entirely.
#
# Sheets "fundamental_data" and "portfolio_data":
#   - JSW, Grupa Kety, Eurocash and CitiBank Handlowy are removed entirely.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: target_data -------------------------------------------------
$wsTarget = $wb.Worksheets.Item("target_data")

# Delete rows bottom-to-top so earlier row numbers stay valid.
$wsTarget.Rows.Item(29).Delete()   # CitiBank Handlowy
$wsTarget.Rows.Item(28).Delete()   # Eurocash
$wsTarget.Rows.Item(27).Delete()   # Grupa Kety
$wsTarget.Rows.Item(26).Delete()   # JSW
$wsTarget.Rows.Item(9).Delete()    # Orange, scope S1+S3 (duplicate row)
$wsTarget.Rows.Item(3).Delete()    # Allegro, scope S3 (duplicate row)

# --- Sheet 2: fundamental_data ---------------------------------------------
$wsFund = $wb.Worksheets.Item("fundamental_data")

$wsFund.Rows.Item(25).Delete()     # CitiBank Handlowy
$wsFund.Rows.Item(24).Delete()     # Eurocash
$wsFund.Rows.Item(23).Delete()     # Grupa Kety
$wsFund.Rows.Item(22).Delete()     # JSW

# --- Sheet 3: portfolio_data ------------------------------------------------
$wsPort = $wb.Worksheets.Item("portfolio_data")

$wsPort.Rows.Item(25).Delete()     # CitiBank Handlowy
$wsPort.Rows.Item(24).Delete()     # Eurocash
$wsPort.Rows.Item(23).Delete()     # Grupa Kety
$wsPort.Rows.Item(22).Delete()     # JSW
